$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "greedy_runtime"

# Extend formatted style (border/bold/centered) from A10 down through A20
$ws.Range("A10").Copy($ws.Range("A11:A20"))

# Update existing rows 2-10: column C becomes a text "time   throughput" value
$ws.Range("C2").Value = "0.000431060791015625   1478.7844496356363"
$ws.Range("C3").Value = "0.0017931461334228516   1586.6035551964073"
$ws.Range("C4").Value = "0.0010921955108642578   2700.1695782419783"
$ws.Range("C5").Value = "0.0009570121765136719   2304.1943061377356"
$ws.Range("C6").Value = "0.0009441375732421875   2456.983204401306"
$ws.Range("C7").Value = "0.001055002212524414   3115.6124635030487"
$ws.Range("C8").Value = "0.0008959770202636719   3506.6193497047925"
$ws.Range("C9").Value = "0.001039743423461914   3058.6771489783328"
$ws.Range("C10").Value = "0.0011591911315917969   2949.182516313036"

# Add new rows 11-20
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = "0.0007276535034179688   2614.7657636059994"
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 13
$ws.Range("C12").Value = "0.0006642341613769531   3440.302738668852"
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 14
$ws.Range("C13").Value = "0.0023970603942871094   4064.985423171131"
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = "0.0012159347534179688   4610.146638695583"
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 16
$ws.Range("C15").Value = "0.0018970966339111328   4460.684505594079"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 32
$ws.Range("C16").Value = "0.0012981891632080078   6986.446840855578"
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 64
$ws.Range("C17").Value = "0.002807140350341797   8367.502439442946"
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 128
$ws.Range("C18").Value = "0.007539987564086914   9836.614937019467"
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 256
$ws.Range("C19").Value = "0.012941122055053711   15922.881833623529"
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 512
$ws.Range("C20").Value = "0.05826616287231445   21238.789086831857"
